$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.489.34"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "2.531.25"
$ws.Range("E3").Value = "  -1.58%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'574.94"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").Value = "'168.53"
$ws.Range("E6").Value = "  -1.60%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.522"
$ws.Range("E8").Value = "  +1.97%  "
$ws.Range("D9").Value = "2.534.03"
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("E10").Value = "  -2.29%  "
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").Value = "'0.357"
$ws.Range("E12").Value = "  +1.96%  "
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("D14").Value = "3.003.61"
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("D15").Value = "70.574.40"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("E16").Value = "  -3.05%  "
$ws.Range("D17").Value = "'25.25"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "2.543.55"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").Value = "'7.92"
$ws.Range("E19").Value = "  +5.67%  "
$ws.Range("D20").Value = "'11.47"
$ws.Range("E20").Value = "  -3.09%  "
$ws.Range("D21").Value = "'352.01"
$ws.Range("E21").Value = "  -3.48%  "
$ws.Range("D22").Value = "'3.94"
$ws.Range("E22").Value = "  -1.62%  "
$ws.Range("D23").Value = "'2.01"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "'70.41"
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").Value = "'4.03"
$ws.Range("E26").Value = "  -2.95%  "
$ws.Range("D27").Value = "'8.92"
$ws.Range("E27").Value = "  -4.82%  "
$ws.Range("D28").Value = "2.686.41"
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").Value = "0.0₃0912"
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("D31").Value = "'7.92"
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("D32").Value = "'466.31"
$ws.Range("E32").Value = "  -4.34%  "
$ws.Range("D33").Value = "'1.26"
$ws.Range("E33").Value = "  -3.89%  "
$ws.Range("D34").Value = "'1.75"
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("E36").Value = "  +3.42%  "
$ws.Range("D37").Value = "'158.12"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D38").Value = "'19.09"
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").Value = "'18.78"
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "'4.81"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("D42").Value = "'0.319"
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("D43").Value = "'1.61"
$ws.Range("E43").Value = "  -5.04%  "
$ws.Range("D44").Value = "'2.33"
$ws.Range("E44").Value = "  -6.67%  "
$ws.Range("E45").Value = "  -13.50%  "
$ws.Range("D46").Value = "'38.55"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "'143.40"
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("D48").Value = "'0.531"
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("D49").Value = "'3.52"
$ws.Range("E49").Value = "  -2.07%  "
$ws.Range("D50").Value = "'1.60"
$ws.Range("E50").Value = "  -3.52%  "
$ws.Range("D51").Value = "'0.0734"
$ws.Range("E51").Value = "  -0.33%  "
